$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B4").Value = 167.94
$ws.Range("C4").Value = 4597.3999999999996
